# Updated cryptos list - applies diff changes to D/E columns (and B/C for a row reorder)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.672.04"
$ws.Range('E2').Value = "'  +0.70%  "
$ws.Range('D3').Value = "'1.805.88"
$ws.Range('E3').Value = "'  -0.34%  "
$ws.Range('E4').Value = "'  +0.16%  "
$ws.Range('D5').Value = "'317.87"
$ws.Range('E5').Value = "'  +0.31%  "
$ws.Range('E6').Value = "'  +0.09%  "
$ws.Range('D7').Value = "'0.5416"
$ws.Range('E7').Value = "'  -2.21%  "
$ws.Range('D8').Value = "'0.3792"
$ws.Range('E8').Value = "'  -1.75%  "
$ws.Range('D9').Value = "'0.07513"
$ws.Range('E9').Value = "'  -1.17%  "
$ws.Range('D10').Value = "'42.44"
$ws.Range('E10').Value = "'  -1.15%  "
$ws.Range('D11').Value = "'1.114"
$ws.Range('E11').Value = "'  -2.01%  "
$ws.Range('D12').Value = "'1.001"
$ws.Range('E12').Value = "'  +0.15%  "
$ws.Range('D13').Value = "'20.65"
$ws.Range('E13').Value = "'  -2.75%  "
$ws.Range('D14').Value = "'6.156"
$ws.Range('E14').Value = "'  -1.47%  "
$ws.Range('D15').Value = "'7.362"
$ws.Range('E15').Value = "'  +0.49%  "
$ws.Range('D16').Value = "'1.799.54"
$ws.Range('E16').Value = "'  -0.65%  "
$ws.Range('D17').Value = "'90.19"
$ws.Range('E17').Value = "'  -1.28%  "
$ws.Range('D18').Value = "'0.00001066"
$ws.Range('E18').Value = "'  -0.84%  "
$ws.Range('D19').Value = "'0.06505"
$ws.Range('E19').Value = "'  +0.38%  "
$ws.Range('E20').Value = "'  +0.03%  "
$ws.Range('D21').Value = "'17.37"
$ws.Range('E21').Value = "'  +0.48%  "
$ws.Range('D23').Value = "'28.678.56"
$ws.Range('E23').Value = "'  +0.77%  "
$ws.Range('E24').Value = "'  -1.58%  "
$ws.Range('D25').Value = "'2.075"
$ws.Range('E25').Value = "'  -2.53%  "
$ws.Range('D26').Value = "'161.21"
$ws.Range('E26').Value = "'  +2.52%  "
$ws.Range('D27').Value = "'20.50"
$ws.Range('E27').Value = "'  -1.18%  "
$ws.Range('D28').Value = "'2.007.83"
$ws.Range('E28').Value = "'  -0.64%  "
$ws.Range('D29').Value = "'2.338"
$ws.Range('E29').Value = "'  -4.33%  "
$ws.Range('D30').Value = "'123.15"
$ws.Range('E30').Value = "'  -0.83%  "
$ws.Range('D31').Value = "'1.150"
$ws.Range('E31').Value = "'  -1.77%  "
$ws.Range('D32').Value = "'0.1055"
$ws.Range('E32').Value = "'  +1.82%  "
$ws.Range('D33').Value = "'5.665"
$ws.Range('E33').Value = "'  -2.00%  "
$ws.Range('D34').Value = "'3.684"
$ws.Range('E34').Value = "'  +1.15%  "
$ws.Range('D35').Value = "'0.06657"
$ws.Range('E35').Value = "'  +6.93%  "
$ws.Range('D36').Value = "'0.2271"
$ws.Range('E36').Value = "'  -0.58%  "
$ws.Range('D37').Value = "'0.02311"
$ws.Range('E37').Value = "'  -0.99%  "
$ws.Range('D38').Value = "'8.622"
$ws.Range('E38').Value = "'  -3.32%  "
$ws.Range('D39').Value = "'5.034"
$ws.Range('E39').Value = "'  -0.13%  "
$ws.Range('B40').Value = "'TheSandbox"
$ws.Range('C40').Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range('D40').Value = "'0.6238"
$ws.Range('E40').Value = "'  -2.57%  "
$ws.Range('B41').Value = "'Aptos"
$ws.Range('C41').Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('D41').Value = "'11.25"
$ws.Range('E41').Value = "'  -3.52%  "
$ws.Range('D42').Value = "'1.201"
$ws.Range('E42').Value = "'  +1.94%  "
$ws.Range('D43').Value = "'1.451"
$ws.Range('E43').Value = "'  +5.01%  "
$ws.Range('D44').Value = "'1.000"
$ws.Range('E44').Value = "'  +0.01%  "
$ws.Range('D45').Value = "'13.34"
$ws.Range('E45').Value = "'  -0.41%  "
$ws.Range('D46').Value = "'3.704"
$ws.Range('E46').Value = "'  +0.10%  "
$ws.Range('D47').Value = "'0.5856"
$ws.Range('E47').Value = "'  -2.54%  "
$ws.Range('D48').Value = "'126.75"
$ws.Range('E48').Value = "'  +2.85%  "
$ws.Range('D49').Value = "'1.960"
$ws.Range('E49').Value = "'  -0.67%  "
$ws.Range('D50').Value = "'1.187"
$ws.Range('E50').Value = "'  +3.59%  "
$ws.Range('E51').Value = "'  -0.31%  "
